$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.739.07'
$ws.Range("E2").Value = '  +0.78%  '

$ws.Range("D3").Value = '3.426.92'
$ws.Range("E3").Value = '  +1.19%  '

$ws.Range("E4").Value = '  +0.13%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '583.84'
$ws.Range("E5").Value = '  -0.51%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '176.94'
$ws.Range("E6").Value = '  -1.79%  '

$ws.Range("D7").Value = '3.419.56'
$ws.Range("E7").Value = '  +1.15%  '

$ws.Range("E8").Value = '  +0.07%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.597'
$ws.Range("E9").Value = '  +0.24%  '

$ws.Range("E10").Value = '  +3.02%  '

$ws.Range("E11").Value = '  -1.05%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '48.91'
$ws.Range("E12").Value = '  +0.69%  '

$ws.Range("E13").Value = '  +1.02%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '691.63'
$ws.Range("E14").Value = '  +2.19%  '

$ws.Range("D15").Value = '3.973.48'
$ws.Range("E15").Value = '  +0.99%  '

$ws.Range("E16").Value = '  +0.25%  '

$ws.Range("D17").Value = '69.747.41'
$ws.Range("E17").Value = '  +0.72%  '

$ws.Range("D18").Value = '3.430.55'
$ws.Range("E18").Value = '  +1.11%  '

$ws.Range("E19").Value = '  +1.09%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.68'
$ws.Range("E20").Value = '  +0.10%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.42'
$ws.Range("E21").Value = '  +0.67%  '

$ws.Range("E22").Value = '  -0.29%  '

$ws.Range("E23").Value = '  +0.90%  '

$ws.Range("E24").Value = '  -0.85%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '100.80'
$ws.Range("E25").Value = '  -2.50%  '

$ws.Range("E26").Value = '  +0.15%  '

$ws.Range("E27").Value = '  -2.28%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.61'
$ws.Range("E28").Value = '  +0.26%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '33.49'
$ws.Range("E29").Value = '  -1.78%  '

$ws.Range("E30").Value = '  +0.77%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.12'
$ws.Range("E31").Value = '  +1.92%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '574.62'
$ws.Range("E32").Value = '  +3.75%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.74'
$ws.Range("E33").Value = '  +0.64%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '11.01'
$ws.Range("E34").Value = '  -1.53%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '58.31'
$ws.Range("E35").Value = '  +0.59%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.103'
$ws.Range("E36").Value = '  -2.45%  '

$ws.Range("D38").Value = '3.589.87'
$ws.Range("E38").Value = '  -2.78%  '

$ws.Range("E39").Value = '  +0.01%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '35.21'
$ws.Range("E40").Value = '  +0.27%  '

$ws.Range("E41").Value = '  +4.77%  '

$ws.Range("E42").Value = '  +0.77%  '

$ws.Range("E43").Value = '  +0.33%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.38'
$ws.Range("E44").Value = '  +3.38%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.333'
$ws.Range("E45").Value = '  -1.32%  '

$ws.Range("E46").Value = '  -0.52%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.47'
$ws.Range("E47").Value = '  +4.66%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.67'
$ws.Range("E48").Value = '  +0.53%  '

$ws.Range("E49").Value = '  -0.31%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.00'
$ws.Range("E50").Value = '  -0.14%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '133.06'
$ws.Range("E51").Value = '  +0.98%  '
